# TEMPLATE 3 - BILAN SENOLOGIQUE SUSPECT
#
# The template used to be a single hard-coded sample report (one
# specific patient, one specific finding). This rewrites it into a
# generic, pipeline-fillable template: a short title, then one
# paragraph per report section (Indication, Mammographie, Echographie,
# Conclusion...) with "[placeholder]" fields a pipeline can substitute.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Wipe the existing body content.
#
# Word never lets a document end up with zero paragraphs, so the
# final paragraph mark survives every Content.Delete() as an
# (eventually empty) trailing paragraph. Two passes collapse the two
# pre-existing paragraphs down to that single empty one, which is
# exactly the bare trailing <w:p/> the new layout ends with, so it is
# left in place rather than removed.
# ---------------------------------------------------------------------
$d.Content.Delete()
$d.Content.Delete()

# ---------------------------------------------------------------------
# 2) Rebuild the report as WordprocessingML and drop it in at the top
#    via Range.InsertXML(), which lets us set exact run/paragraph
#    formatting (bold/underline title, the superscript "e" in "2e"
#    rendered with an Arial run, the grammar/spelling proofErr spans,
#    xml:space="preserve" runs, line breaks, etc.) in one shot.
# ---------------------------------------------------------------------

$paraTitle = '<w:p>' +
    '<w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr>' +
    '<w:t xml:space="preserve">BILAN SÉNOLOGIQUE </w:t></w:r>' +
    '</w:p>'

$paraSubtitle = '<w:p>' +
    '<w:r><w:t>Compte rendu avant 2</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>ᵉ</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> lecture du dépistage organisé</w:t></w:r>' +
    '</w:p>'

$paraIndication = '<w:p>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Indication :</w:t></w:r>' +
    '<w:r><w:br/><w:t xml:space="preserve">Patiente de [] </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>ans adressée</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> pour [motif].</w:t></w:r>' +
    '<w:r><w:br/><w:t>Parité : [] Allaitement : [Oui / Non]</w:t></w:r>' +
    '<w:r><w:br/><w:t>ATCD personnels : [ATCD pers]</w:t></w:r>' +
    '<w:r><w:br/><w:t xml:space="preserve">ATCD familiaux : [ATCD </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>fam</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>]</w:t></w:r>' +
    '<w:r><w:br/><w:t>Examen clinique : [Examen]</w:t></w:r>' +
    '</w:p>'

$paraMammo = '<w:p>' +
    '<w:r><w:t>Mammographie bilatérale numérisée</w:t></w:r>' +
    '<w:r><w:br/><w:t>Densité mammaire ACR : [ACR densité]</w:t></w:r>' +
    '<w:r><w:br/><w:t>Sein droit : [Description mammo droite]</w:t></w:r>' +
    '<w:r><w:br/><w:t>Sein gauche : [Description mammo gauche]</w:t></w:r>' +
    '</w:p>'

$paraEcho = '<w:p>' +
    '<w:r><w:t>Échographie mammaire</w:t></w:r>' +
    '<w:r><w:br/><w:t>Sein droit : [Description écho droite]</w:t></w:r>' +
    '<w:r><w:br/><w:t>Sein gauche : [Description écho gauche]</w:t></w:r>' +
    '<w:r><w:br/><w:t>Axillaires : [Description]</w:t></w:r>' +
    '</w:p>'

$paraConclusionHeading = '<w:p>' +
    '<w:pPr><w:spacing w:after="0"/></w:pPr>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>CONCLUSION :</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '</w:p>'

$paraConclusionBody = '<w:p>' +
    '<w:r><w:t>[Conclusion détaillée]</w:t></w:r>' +
    '<w:r><w:br/><w:t>ACR [] à droite.</w:t></w:r>' +
    '<w:r><w:br/><w:t>ACR [] à gauche.</w:t></w:r>' +
    '<w:r><w:br/><w:t>Recommandation : [biopsie / contrôle / autre]</w:t></w:r>' +
    '</w:p>'

$bodyParagraphs = $paraTitle + $paraSubtitle + $paraIndication + $paraMammo +
    $paraEcho + $paraConclusionHeading + $paraConclusionBody

$newBodyXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $bodyParagraphs + '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint = $d.Range(0, 0)
[void]$insertionPoint.InsertXML($newBodyXml)
